# Fix the mislabeled column-E header ("2050" / "2041-2050") on the
# scenario tables and remove the "Total" summary rows that should no
# longer be present.
$wb = $excel.ActiveWorkbook

# Sheets whose column-E (row 1) header text must be corrected, and whether
# they also carry a "Total" row (row 13) that must be removed.
$sheetInfo = @{
    "Potencia Acumulada - SIN (MW)"   = @{ Label = "2050";      HasTotalRow13 = $true }
    "Geracao Periodo Medio (MWMed)"   = @{ Label = "2050";      HasTotalRow13 = $true }
    "Atendimento a Ponta(MW)"         = @{ Label = "2050";      HasTotalRow13 = $true }
    "Potencia Incremental - SIN(MW)"  = @{ Label = "2041-2050"; HasTotalRow13 = $true }
    "Emissoes Totais (MtCO2eq)"       = @{ Label = "2050";      HasTotalRow13 = $false }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name

    if ($sheetInfo.ContainsKey($name)) {
        $info = $sheetInfo[$name]

        # Write the corrected label as TEXT (the other headers in row 1,
        # e.g. "2015"/"2030"/"2040", are text too) while keeping the
        # cell's existing style (bold, centered, bordered - style index
        # carried by D1, its neighbour). A plain .Value assignment would
        # get auto-coerced to a number since the label is numeric-looking,
        # so force Text format first, then restore the original look via
        # a formats-only paste from the untouched neighbour cell.
        $ws.Range("E1").NumberFormat = "@"
        $ws.Range("E1").Value = $info.Label
        $ws.Range("D1").Copy()
        $ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
        $excel.CutCopyMode = $false

        if ($info.HasTotalRow13) {
            $ws.Rows.Item(13).Delete()
        }
    } elseif ($name -eq "Custo Total (bilhões de R$)") {
        # This sheet's total row is row 4 (columns A:B)
        $ws.Rows.Item(4).Delete()
    }
}
